$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(12, 8).Value = 322  # H12
$ws.Cells.Item(12, 9).Value = 322  # I12
$ws.Cells.Item(12, 10).Value = 0  # J12
$ws.Cells.Item(12, 11).Value = 322  # K12
$ws.Cells.Item(12, 12).Value = 0  # L12
$ws.Cells.Item(12, 13).Value = -152  # M12
$ws.Cells.Item(12, 14).Value = $null  # N12
$ws.Cells.Item(114, 8).Value = 0  # H114
$ws.Cells.Item(114, 10).Value = 0  # J114
$ws.Cells.Item(114, 12).Value = 0  # L114
$ws.Cells.Item(114, 14).Value = $null  # N114
$ws.Cells.Item(116, 8).Value = 2150  # H116
$ws.Cells.Item(116, 9).Value = 2133.3333  # I116
$ws.Cells.Item(116, 11).Value = 2133.3333  # K116
$ws.Cells.Item(116, 13).Value = 1308.6667  # M116
$ws.Cells.Item(117, 8).Value = 0  # H117
$ws.Cells.Item(117, 10).Value = 0  # J117
$ws.Cells.Item(117, 12).Value = 0  # L117
$ws.Cells.Item(117, 14).Value = $null  # N117
$ws.Cells.Item(129, 8).Value = 17826.271  # H129
$ws.Cells.Item(129, 9).Value = 461.4  # I129
$ws.Cells.Item(129, 10).Value = 21370.123  # J129
$ws.Cells.Item(129, 11).Value = 1384.2  # K129
$ws.Cells.Item(129, 12).Value = 64110.369  # L129
$ws.Cells.Item(129, 13).Value = 3615.8  # M129
$ws.Cells.Item(129, 14).Value = -74110.36900000001  # N129
$ws.Cells.Item(137, 8).Value = 2015.4166  # H137
$ws.Cells.Item(137, 9).Value = 1998.2  # I137
$ws.Cells.Item(137, 10).Value = 2101.5  # J137
$ws.Cells.Item(137, 11).Value = 5994.6  # K137
$ws.Cells.Item(137, 12).Value = 6304.5  # L137
$ws.Cells.Item(137, 13).Value = -3444.6  # M137
$ws.Cells.Item(137, 14).Value = -11404.5  # N137

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 16334.931  # H32
$ws.Cells.Item(32, 9).Value = 16871.232  # I32
$ws.Cells.Item(32, 10).Value = 4000  # J32
$ws.Cells.Item(32, 11).Value = 16871.232  # K32
$ws.Cells.Item(32, 12).Value = 4000  # L32
$ws.Cells.Item(32, 13).Value = -16584.232  # M32
$ws.Cells.Item(32, 14).Value = -4574  # N32
$ws.Cells.Item(61, 8).Value = 2040  # H61
$ws.Cells.Item(61, 9).Value = 1074.2858  # I61
$ws.Cells.Item(61, 10).Value = 2654.5454  # J61
$ws.Cells.Item(61, 11).Value = 1074.2858  # K61
$ws.Cells.Item(61, 12).Value = 2654.5454  # L61
$ws.Cells.Item(61, 13).Value = -862.2858000000001  # M61
$ws.Cells.Item(61, 14).Value = -3078.5454  # N61
$ws.Cells.Item(74, 8).Value = 1432.7693  # H74
$ws.Cells.Item(74, 9).Value = 1628.25  # I74
$ws.Cells.Item(74, 10).Value = 1120  # J74
$ws.Cells.Item(74, 11).Value = 1628.25  # K74
$ws.Cells.Item(74, 12).Value = 1120  # L74
$ws.Cells.Item(74, 13).Value = -754.25  # M74
$ws.Cells.Item(74, 14).Value = -2868  # N74
$ws.Cells.Item(77, 8).Value = 1432.7693  # H77
$ws.Cells.Item(77, 9).Value = 1628.25  # I77
$ws.Cells.Item(77, 10).Value = 1120  # J77
$ws.Cells.Item(77, 11).Value = 8141.25  # K77
$ws.Cells.Item(77, 12).Value = 5600  # L77
$ws.Cells.Item(77, 13).Value = -3773.25  # M77
$ws.Cells.Item(77, 14).Value = -14336  # N77
$ws.Cells.Item(108, 8).Value = 30000  # H108
$ws.Cells.Item(108, 10).Value = 30000  # J108
$ws.Cells.Item(108, 12).Value = 30000  # L108
$ws.Cells.Item(108, 14).Value = -37680  # N108
$ws.Cells.Item(136, 8).Value = 2040  # H136
$ws.Cells.Item(136, 9).Value = 1074.2858  # I136
$ws.Cells.Item(136, 10).Value = 2654.5454  # J136
$ws.Cells.Item(136, 11).Value = 3222.8574  # K136
$ws.Cells.Item(136, 12).Value = 7963.6362  # L136
$ws.Cells.Item(136, 13).Value = -672.8574000000003  # M136
$ws.Cells.Item(136, 14).Value = -13063.6362  # N136

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(122, 8).Value = 430  # H122
$ws.Cells.Item(122, 9).Value = 452  # I122
$ws.Cells.Item(122, 10).Value = 298  # J122
$ws.Cells.Item(122, 11).Value = 1356  # K122
$ws.Cells.Item(122, 12).Value = 894  # L122
$ws.Cells.Item(122, 13).Value = 1094  # M122
$ws.Cells.Item(122, 14).Value = -5794  # N122

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 5556749.5  # H131
$ws.Cells.Item(131, 9).Value = 1772.5  # I131
$ws.Cells.Item(131, 10).Value = 9260068  # J131
$ws.Cells.Item(131, 11).Value = 5317.5  # K131
$ws.Cells.Item(131, 12).Value = 27780204  # L131
$ws.Cells.Item(131, 13).Value = -277.5  # M131
$ws.Cells.Item(131, 14).Value = -27790284  # N131
$ws.Cells.Item(137, 8).Value = 55568748  # H137
$ws.Cells.Item(137, 9).Value = 33334676  # I137
$ws.Cells.Item(137, 10).Value = 70391464  # J137
$ws.Cells.Item(137, 11).Value = 100004028  # K137
$ws.Cells.Item(137, 12).Value = 211174392  # L137
$ws.Cells.Item(137, 13).Value = -99998928  # M137
$ws.Cells.Item(137, 14).Value = -211184592  # N137

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(32, 8).Value = 0  # H32
$ws.Cells.Item(32, 10).Value = 0  # J32
$ws.Cells.Item(32, 12).Value = 0  # L32
$ws.Cells.Item(32, 14).Value = $null  # N32
$ws.Cells.Item(42, 8).Value = 34851.477  # H42
$ws.Cells.Item(42, 10).Value = 34851.477  # J42
$ws.Cells.Item(42, 12).Value = 34851.477  # L42
$ws.Cells.Item(42, 14).Value = -35821.477  # N42
$ws.Cells.Item(113, 8).Value = 16669301  # H113
$ws.Cells.Item(113, 9).Value = 31252376  # I113
$ws.Cells.Item(113, 10).Value = 2928.5715  # J113
$ws.Cells.Item(113, 11).Value = 31252376  # K113
$ws.Cells.Item(113, 12).Value = 2928.5715  # L113
$ws.Cells.Item(113, 13).Value = -31250206  # M113
$ws.Cells.Item(113, 14).Value = -7268.5715  # N113
$ws.Cells.Item(115, 8).Value = 34851.477  # H115
$ws.Cells.Item(115, 10).Value = 34851.477  # J115
$ws.Cells.Item(115, 12).Value = 34851.477  # L115
$ws.Cells.Item(115, 14).Value = -37201.477  # N115
$ws.Cells.Item(126, 8).Value = 1099.4445  # H126
$ws.Cells.Item(126, 10).Value = 1379  # J126
$ws.Cells.Item(126, 12).Value = 4137  # L126
$ws.Cells.Item(126, 14).Value = -9077  # N126

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 2498.182  # H7
$ws.Cells.Item(7, 9).Value = 2100  # I7
$ws.Cells.Item(7, 10).Value = 2725.7144  # J7
$ws.Cells.Item(7, 11).Value = 2100  # K7
$ws.Cells.Item(7, 12).Value = 2725.7144  # L7
$ws.Cells.Item(7, 13).Value = -1988  # M7
$ws.Cells.Item(7, 14).Value = -2949.7144  # N7
$ws.Cells.Item(61, 8).Value = 1152.5883  # H61
$ws.Cells.Item(61, 9).Value = 1199.9231  # I61
$ws.Cells.Item(61, 10).Value = 998.75  # J61
$ws.Cells.Item(61, 11).Value = 1199.9231  # K61
$ws.Cells.Item(61, 12).Value = 998.75  # L61
$ws.Cells.Item(61, 13).Value = -997.9231  # M61
$ws.Cells.Item(61, 14).Value = -1402.75  # N61
$ws.Cells.Item(113, 8).Value = 1152.5883  # H113
$ws.Cells.Item(113, 9).Value = 1199.9231  # I113
$ws.Cells.Item(113, 10).Value = 998.75  # J113
$ws.Cells.Item(113, 11).Value = 1199.9231  # K113
$ws.Cells.Item(113, 12).Value = 998.75  # L113
$ws.Cells.Item(113, 13).Value = 970.0769  # M113
$ws.Cells.Item(113, 14).Value = -5338.75  # N113
$ws.Cells.Item(122, 8).Value = 3614.3157  # H122
$ws.Cells.Item(122, 9).Value = 4667.7144  # I122
$ws.Cells.Item(122, 10).Value = 2999.8333  # J122
$ws.Cells.Item(122, 11).Value = 14003.1432  # K122
$ws.Cells.Item(122, 12).Value = 8999.499899999999  # L122
$ws.Cells.Item(122, 13).Value = -11553.1432  # M122
$ws.Cells.Item(122, 14).Value = -13899.4999  # N122
$ws.Cells.Item(126, 8).Value = 2498.182  # H126
$ws.Cells.Item(126, 9).Value = 2100  # I126
$ws.Cells.Item(126, 10).Value = 2725.7144  # J126
$ws.Cells.Item(126, 11).Value = 6300  # K126
$ws.Cells.Item(126, 12).Value = 8177.1432  # L126
$ws.Cells.Item(126, 13).Value = -3830  # M126
$ws.Cells.Item(126, 14).Value = -13117.1432  # N126
$ws.Cells.Item(136, 8).Value = 8986.933999999999  # H136
$ws.Cells.Item(136, 9).Value = 13111.556  # I136
$ws.Cells.Item(136, 10).Value = 2800  # J136
$ws.Cells.Item(136, 11).Value = 39334.66800000001  # K136
$ws.Cells.Item(136, 12).Value = 8400  # L136
$ws.Cells.Item(136, 13).Value = -36784.66800000001  # M136
$ws.Cells.Item(136, 14).Value = -13500  # N136
$ws.Cells.Item(139, 8).Value = 42085.75  # H139
$ws.Cells.Item(139, 9).Value = 39000  # I139
$ws.Cells.Item(139, 10).Value = 42248.156  # J139
$ws.Cells.Item(139, 11).Value = 39000  # K139
$ws.Cells.Item(139, 12).Value = 42248.156  # L139
$ws.Cells.Item(139, 13).Value = -33860  # M139
$ws.Cells.Item(139, 14).Value = -52528.156  # N139

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(27, 8).Value = 54582  # H27
$ws.Cells.Item(27, 9).Value = 0  # I27
$ws.Cells.Item(27, 10).Value = 54582  # J27
$ws.Cells.Item(27, 11).Value = 0  # K27
$ws.Cells.Item(27, 12).Value = 54582  # L27
$ws.Cells.Item(27, 13).Value = $null  # M27
$ws.Cells.Item(27, 14).Value = -54720  # N27
$ws.Cells.Item(113, 8).Value = 403.92307  # H113
$ws.Cells.Item(113, 9).Value = 429.25  # I113
$ws.Cells.Item(113, 10).Value = 100  # J113
$ws.Cells.Item(113, 11).Value = 1287.75  # K113
$ws.Cells.Item(113, 12).Value = 300  # L113
$ws.Cells.Item(113, 13).Value = 882.25  # M113
$ws.Cells.Item(113, 14).Value = -4640  # N113
$ws.Cells.Item(115, 8).Value = 38604  # H115
$ws.Cells.Item(115, 10).Value = 38604  # J115
$ws.Cells.Item(115, 12).Value = 38604  # L115
$ws.Cells.Item(115, 14).Value = -41738  # N115
$ws.Cells.Item(126, 8).Value = 13163  # H126
$ws.Cells.Item(126, 9).Value = 14900.571  # I126
$ws.Cells.Item(126, 10).Value = 1000  # J126
$ws.Cells.Item(126, 11).Value = 44701.713  # K126
$ws.Cells.Item(126, 12).Value = 3000  # L126
$ws.Cells.Item(126, 13).Value = -42231.713  # M126
$ws.Cells.Item(126, 14).Value = -7940  # N126
$ws.Cells.Item(136, 8).Value = 9869.857  # H136
$ws.Cells.Item(136, 9).Value = 11398.167  # I136
$ws.Cells.Item(136, 11).Value = 34194.501  # K136
$ws.Cells.Item(136, 13).Value = -31644.501  # M136
